$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The picture is anchored below the insertion point, so it must shift down
# by the height of the newly inserted rows (4 rows). Capture its current
# Top plus a row-height reference before the insert so we can compute the
# exact shift afterwards.
$shp = $ws.Shapes.Item(1)
$origTop = $shp.Top
$row10TopBefore = $ws.Range("A10").Top

# Insert 4 new rows before row 5; existing rows 5-8 shift down to rows
# 9-12 (taking their values/styles with them), and everything below -
# including the picture - moves down accordingly.
$ws.Rows("5:8").Insert()

$row14TopAfter = $ws.Range("A14").Top
$rowShift = $row14TopAfter - $row10TopBefore
$shp.Top = $origTop + $rowShift

# The newly inserted rows copy the formatting of the row above (row 4);
# clear that so the new rows start out unformatted/blank.
$ws.Range("A5:J8").Clear()

# Populate the new "Select <Tier>" rows that were inserted above the
# existing "Choose <Tier>" rows (now rows 9-12).
$ws.Range("A5").Value = "Select Silver"
$ws.Range("C5").Value = "<CHECK>"

$ws.Range("A6").Value = "Select Gold"
$ws.Range("D6").Value = "<CHECK>"

$ws.Range("A7").Value = "Select Platinum"
$ws.Range("E7").Value = "<CHECK>"

$ws.Range("A8").Value = "Select Ultimate"
$ws.Range("F8").Value = "<CHECK>"

# Update the active selection to match the edited workbook.
[void]$ws.Range("A5:A8").Select()

Write-Host "Done"
